# For every data row in the active sheet, the "Recorded By" column (G) holds
# a comma-separated list of names/emails (e.g. "dnasr281@gmail.com, System").
# This edit rotates each such list one position to the right - i.e. the last
# entry in the list is moved to the front - leaving single-entry cells
# untouched. This mirrors the author's reordering of the "Recorded By"
# attribution lists across the whole "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ', '

    if ($parts.Count -gt 1) {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $newValue = [string]::Join(', ', $rotated)
        $cell.Value = $newValue
    }
}
